# Add data for 2024-12-24
# Updates 2024 (K) year-to-date figures (and a few other prior columns)
# across the Citywide Totals, By Neighborhood summary, and individual
# neighborhood sheets, to reflect the newly added day of data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7708
$ws.Range("K3").Value = 7968
$ws.Range("E4").Value = 2043
$ws.Range("J4").Value = 1849
$ws.Range("K4").Value = 1672
$ws.Range("J6").Value = 11054
$ws.Range("K6").Value = 8877
$ws.Range("E7").Value = 26049
$ws.Range("J7").Value = 29319
$ws.Range("K7").Value = 26795

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K8").Value = 1750
$ws.Range("K17").Value = 51
$ws.Range("K18").Value = 182
$ws.Range("K19").Value = 771
$ws.Range("K20").Value = 659
$ws.Range("K21").Value = 91
$ws.Range("K23").Value = 265
$ws.Range("K25").Value = 126
$ws.Range("K29").Value = 1482
$ws.Range("K33").Value = 1125
$ws.Range("K36").Value = 347
$ws.Range("K37").Value = 879
$ws.Range("K44").Value = 215
$ws.Range("K48").Value = 333
$ws.Range("K52").Value = 691
$ws.Range("K54").Value = 527
$ws.Range("K55").Value = 294
$ws.Range("E63").Value = 377
$ws.Range("J63").Value = 196
$ws.Range("K64").Value = 160
$ws.Range("K66").Value = 79
$ws.Range("K67").Value = 1040
$ws.Range("K70").Value = 48
$ws.Range("K72").Value = 126
$ws.Range("K73").Value = 240
$ws.Range("K76").Value = 370
$ws.Range("K77").Value = 175
$ws.Range("K79").Value = 656
$ws.Range("K83").Value = 566
$ws.Range("K84").Value = 218
$ws.Range("K85").Value = 1232
$ws.Range("K88").Value = 284
$ws.Range("K89").Value = 402
$ws.Range("J90").Value = 308
$ws.Range("K94").Value = 362
$ws.Range("K95").Value = 446
$ws.Range("K96").Value = 287
$ws.Range("K99").Value = 450
$ws.Range("E101").Value = 26049
$ws.Range("J101").Value = 29319
$ws.Range("K101").Value = 26795

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 287

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 110
$ws.Range("K3").Value = 123
$ws.Range("K6").Value = 121
$ws.Range("K7").Value = 402

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 406
$ws.Range("K6").Value = 303
$ws.Range("K7").Value = 1232

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 191
$ws.Range("K6").Value = 246
$ws.Range("K7").Value = 691

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 488
$ws.Range("K3").Value = 530
$ws.Range("K7").Value = 1750

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 197
$ws.Range("K7").Value = 566

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 283
$ws.Range("K3").Value = 397
$ws.Range("K7").Value = 1125

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K4").Value = 19
$ws.Range("K6").Value = 106
$ws.Range("K7").Value = 446

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 250
$ws.Range("K7").Value = 879

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 119
$ws.Range("K7").Value = 450

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 376
$ws.Range("K7").Value = 1040

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 88
$ws.Range("K7").Value = 218

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 285
$ws.Range("K7").Value = 527

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 416
$ws.Range("K3").Value = 521
$ws.Range("K6").Value = 440
$ws.Range("K7").Value = 1482

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K4").Value = 49
$ws.Range("K7").Value = 333

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K6").Value = 258
$ws.Range("K7").Value = 771

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 215

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K3").Value = 72
$ws.Range("K7").Value = 370

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 86
$ws.Range("K3").Value = 83
$ws.Range("K7").Value = 294

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K3").Value = 92
$ws.Range("K7").Value = 265

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 91

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 215
$ws.Range("K7").Value = 656

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 224
$ws.Range("K7").Value = 659

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 182

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 130
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 347

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K4").Value = 30
$ws.Range("K7").Value = 362

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 240

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 89
$ws.Range("K7").Value = 284

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K3").Value = 71
$ws.Range("J4").Value = 12
$ws.Range("K4").Value = 16
$ws.Range("J7").Value = 308

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 73
$ws.Range("K7").Value = 175
